$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update client name, internal ID, and clear status
$ws.Range("C3").Value = "SICOOB MERIDIONAL"
$ws.Range("H3").Value = "ID_A-3370"
$ws.Range("L3").Value = ""

# Row 4: update client name, internal ID, and clear status
$ws.Range("C4").Value = "SICOOB MERIDIONAL"
$ws.Range("H4").Value = "ID_A-3371"
$ws.Range("L4").Value = ""
